# Auto-generated edit script: updates Atomos Profits market-data values
# per the scheduled-runner data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2148.3064
$ws.Range("I15").Value = 2148.3064
$ws.Range("K15").Value = 6444.9192
$ws.Range("M15").Value = -6275.9192
$ws.Range("H40").Value = 62513.77
$ws.Range("I40").Value = 100747.5
$ws.Range("K40").Value = 100747.5
$ws.Range("M40").Value = -100572.5
$ws.Range("H62").Value = 3390.5557
$ws.Range("I62").Value = 3075
$ws.Range("J62").Value = 4021.6667
$ws.Range("K62").Value = 3075
$ws.Range("L62").Value = 4021.6667
$ws.Range("M62").Value = -2451
$ws.Range("N62").Value = -5269.6667
$ws.Range("H65").Value = 3390.5557
$ws.Range("I65").Value = 3075
$ws.Range("J65").Value = 4021.6667
$ws.Range("K65").Value = 15375
$ws.Range("L65").Value = 20108.3335
$ws.Range("M65").Value = -12255
$ws.Range("N65").Value = -26348.3335
$ws.Range("H137").Value = 2909.975
$ws.Range("I137").Value = 3027.862
$ws.Range("J137").Value = 2599.182
$ws.Range("K137").Value = 9083.585999999999
$ws.Range("L137").Value = 7797.545999999999
$ws.Range("M137").Value = -6533.585999999999
$ws.Range("N137").Value = -12897.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1638.3
$ws.Range("I32").Value = 1475.4166
$ws.Range("J32").Value = 2493.4375
$ws.Range("K32").Value = 1475.4166
$ws.Range("L32").Value = 2493.4375
$ws.Range("M32").Value = -1188.4166
$ws.Range("N32").Value = -3067.4375
$ws.Range("H74").Value = 653.9583
$ws.Range("I74").Value = 604.3182
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 604.3182
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = 269.6818
$ws.Range("N74").Value = -2948
$ws.Range("H77").Value = 653.9583
$ws.Range("I77").Value = 604.3182
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 3021.591
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = 1346.409
$ws.Range("N77").Value = -14736
$ws.Range("H132").Value = 1533.8026
$ws.Range("I132").Value = 1010.86884
$ws.Range("J132").Value = 3660.4
$ws.Range("K132").Value = 3032.60652
$ws.Range("L132").Value = 10981.2
$ws.Range("M132").Value = -502.6065199999998
$ws.Range("N132").Value = -16041.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 375.5
$ws.Range("J22").Value = 261
$ws.Range("L22").Value = 261
$ws.Range("N22").Value = -607
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 1944.4193
$ws.Range("I134").Value = 1550.4814
$ws.Range("J134").Value = 4603.5
$ws.Range("K134").Value = 4651.4442
$ws.Range("L134").Value = 13810.5
$ws.Range("M134").Value = -2116.4442
$ws.Range("N134").Value = -18880.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3862.6875
$ws.Range("I31").Value = 2832.6155
$ws.Range("J31").Value = 4567.4736
$ws.Range("K31").Value = 2832.6155
$ws.Range("L31").Value = 4567.4736
$ws.Range("M31").Value = -2537.6155
$ws.Range("N31").Value = -5157.4736
$ws.Range("H34").Value = 3862.6875
$ws.Range("I34").Value = 2832.6155
$ws.Range("J34").Value = 4567.4736
$ws.Range("K34").Value = 2832.6155
$ws.Range("L34").Value = 4567.4736
$ws.Range("M34").Value = -2630.6155
$ws.Range("N34").Value = -4971.4736
$ws.Range("H58").Value = 7044126.5
$ws.Range("I58").Value = 907.5472
$ws.Range("J58").Value = 27782492
$ws.Range("K58").Value = 907.5472
$ws.Range("L58").Value = 27782492
$ws.Range("M58").Value = -704.5472
$ws.Range("N58").Value = -27782898
$ws.Range("H132").Value = 1399
$ws.Range("I132").Value = 1049.3864
$ws.Range("J132").Value = 2253.611
$ws.Range("K132").Value = 3148.1592
$ws.Range("L132").Value = 6760.833
$ws.Range("M132").Value = -618.1592000000001
$ws.Range("N132").Value = -11820.833
$ws.Range("H136").Value = 7044126.5
$ws.Range("I136").Value = 907.5472
$ws.Range("J136").Value = 27782492
$ws.Range("K136").Value = 2722.6416
$ws.Range("L136").Value = 83347476
$ws.Range("M136").Value = -172.6415999999999
$ws.Range("N136").Value = -83352576
$ws.Range("H140").Value = 44800
$ws.Range("J140").Value = 44800
$ws.Range("L140").Value = 44800
$ws.Range("N140").Value = -55160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 11.5
$ws.Range("I2").Value = 9.875
$ws.Range("J2").Value = 12.583333
$ws.Range("K2").Value = 59.25
$ws.Range("L2").Value = 75.49999800000001
$ws.Range("M2").Value = 53.75
$ws.Range("N2").Value = -301.499998
$ws.Range("H34").Value = 9506.923000000001
$ws.Range("I34").Value = 196.66667
$ws.Range("J34").Value = 12300
$ws.Range("K34").Value = 590.00001
$ws.Range("L34").Value = 36900
$ws.Range("M34").Value = -506.00001
$ws.Range("N34").Value = -37068
$ws.Range("H39").Value = 1739
$ws.Range("J39").Value = 1814.8235
$ws.Range("L39").Value = 5444.470499999999
$ws.Range("N39").Value = -6032.470499999999
$ws.Range("H48").Value = 14751
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 14751
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 44253
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -44753
$ws.Range("H55").Value = 2986.087
$ws.Range("J55").Value = 3514.7368
$ws.Range("L55").Value = 10544.2104
$ws.Range("N55").Value = -10898.2104
$ws.Range("H59").Value = 2978
$ws.Range("J59").Value = 2978
$ws.Range("L59").Value = 8934
$ws.Range("N59").Value = -10014
$ws.Range("H64").Value = 1921.3
$ws.Range("J64").Value = 3000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9540
$ws.Range("H67").Value = 1921.3
$ws.Range("J67").Value = 3000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10872
$ws.Range("H109").Value = 640.7778
$ws.Range("I109").Value = 196.16667
$ws.Range("J109").Value = 1530
$ws.Range("K109").Value = 588.50001
$ws.Range("L109").Value = 4590
$ws.Range("M109").Value = 451.49999
$ws.Range("N109").Value = -6670
$ws.Range("H123").Value = 1223.2285
$ws.Range("J123").Value = 3777.6667
$ws.Range("L123").Value = 11333.0001
$ws.Range("N123").Value = -16233.0001
$ws.Range("H132").Value = 3437.5
$ws.Range("J132").Value = 6666.6665
$ws.Range("L132").Value = 59999.9985
$ws.Range("N132").Value = -65059.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1381.25
$ws.Range("I107").Value = 683.3333
$ws.Range("K107").Value = 683.3333
$ws.Range("M107").Value = 1236.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1625.3175
$ws.Range("I132").Value = 929.78723
$ws.Range("J132").Value = 3668.4375
$ws.Range("K132").Value = 2789.36169
$ws.Range("L132").Value = 11005.3125
$ws.Range("M132").Value = -259.3616900000002
$ws.Range("N132").Value = -16065.3125
$ws.Range("H136").Value = 1956.7241
$ws.Range("I136").Value = 1226.0476
$ws.Range("K136").Value = 3678.142800000001
$ws.Range("M136").Value = -1128.142800000001
$ws.Range("H137").Value = 29555
$ws.Range("J137").Value = 29555
$ws.Range("L137").Value = 29555
$ws.Range("N137").Value = -39755

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 44291.715
$ws.Range("J18").Value = 44291.715
$ws.Range("L18").Value = 44291.715
$ws.Range("N18").Value = -44637.715
$ws.Range("H132").Value = 13771.578
$ws.Range("I132").Value = 2760.516
$ws.Range("J132").Value = 38153.215
$ws.Range("K132").Value = 8281.548000000001
$ws.Range("L132").Value = 114459.645
$ws.Range("M132").Value = -5751.548000000001
$ws.Range("N132").Value = -119519.645
$ws.Range("H136").Value = 1348.129
$ws.Range("I136").Value = 902.2381
$ws.Range("J136").Value = 2284.5
$ws.Range("K136").Value = 2706.7143
$ws.Range("L136").Value = 6853.5
$ws.Range("M136").Value = -156.7143000000001
$ws.Range("N136").Value = -11953.5
